$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-17.39779159986714, 1.809310921724846, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714),
    @(-17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, 2.486827604430293, -17.39779159986714, -17.39779159986714),
    @(-17.39779159986714, 2.111877106300822, 2.103892057324111, -17.39779159986714, 3.570372435805767, -17.39779159986714, 1.259390753157166, -17.39779159986714, 2.107347332340634, -17.39779159986714),
    @(-17.39779159986714, 1.817487802552066, -17.39779159986714, -17.39779159986714, -17.39779159986714, 3.048804086991709, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714),
    @(-17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714),
    @(2.590999858450433, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714),
    @(-17.39779159986714, -17.39779159986714, -17.39779159986714, 1.908735190431188, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714),
    @(3.804748967954401, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714),
    @(-17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, 0.8539543019217481, -17.39779159986714, 4.321920158146348),
    @(-17.39779159986714, -17.39779159986714, -17.39779159986714, 2.824013798977296, -17.39779159986714, 2.713314771441542, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714),
    @(-17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714),
    @(-17.39779159986714, -17.39779159986714, -17.39779159986714, 2.51006874362169, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, 1.961771384124566, -17.39779159986714),
    @(-17.39779159986714, -17.39779159986714, 1.197790061024959, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714),
    @(-17.39779159986714, -17.39779159986714, 1.080223522593595, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714),
    @(-17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, 1.677617038517551, -17.39779159986714),
    @(-17.39779159986714, 2.246488953154157, 2.536344123951755, -17.39779159986714, -17.39779159986714, -17.39779159986714, 1.759646626178945, 2.198989619786181, 2.082393118365781, -17.39779159986714),
    @(-17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, -17.39779159986714, 1.987828595966442, 1.583020178636307, 2.12480717117168, -17.39779159986714),
    @(-17.39779159986714, -17.39779159986714, 1.576480134007623, -17.39779159986714, -17.39779159986714, -17.39779159986714, 1.337172654938654, 1.601693736252362, -17.39779159986714, -17.39779159986714),
    @(-17.39779159986714, 0.8990036486486012, 1.327351616368388, -17.39779159986714, 3.021593598477656, -17.39779159986714, 1.969514361546496, 0.9741594162930921, -17.39779159986714, -17.39779159986714),
    @(-17.39779159986714, 1.027100480421685, -17.39779159986714, 1.793829129802085, -17.39779159986714, 2.369152957968805, 1.929982501254223, -17.39779159986714, -17.39779159986714, -17.39779159986714)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($i + 2, $j + 2).Value = $row[$j]
    }
}
